$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.062550666666667
$ws.Range("H2").Value = 9.187652
$ws.Range("I2").Value = 0.06849600470812313
$ws.Range("J2").Value = 0.06849600470812313
$ws.Range("M2").Value = 106.5570066666667
$ws.Range("N2").Value = 319.67102
$ws.Range("O2").Value = 0.956549115742331
$ws.Range("P2").Value = 0.956549115742331
$ws.Range("Q2").Value = 326.3362318050044
$ws.Range("R2").Value = 2937.02608624504
$ws.Range("S2").Value = 0.06551979273543772
$ws.Range("T2").Value = 0.06551979273543772
$ws.Range("G3").Value = 3.062550666666667
$ws.Range("H3").Value = 9.187652
$ws.Range("I3").Value = 0.06849600470812313
$ws.Range("J3").Value = 0.06849600470812313
$ws.Range("N3").Value = 0.7487489999999999
$ws.Range("O3").Value = 0.002240475830004717
$ws.Range("P3").Value = 0.002240475830004717
$ws.Range("Q3").Value = 0.7643605830386666
$ws.Range("R3").Value = 6.879245247347999
$ws.Range("S3").Value = 0.0001534636430004392
$ws.Range("T3").Value = 0.0001534636430004392
$ws.Range("G4").Value = 3.062550666666667
$ws.Range("H4").Value = 9.187652
$ws.Range("I4").Value = 0.06849600470812313
$ws.Range("J4").Value = 0.06849600470812313
$ws.Range("M4").Value = 3.632925333333334
$ws.Range("N4").Value = 10.898776
$ws.Range("O4").Value = 0.03261232296087941
$ws.Range("P4").Value = 0.03261232296087941
$ws.Range("Q4").Value = 11.12601790155022
$ws.Range("R4").Value = 100.134161113952
$ws.Range("S4").Value = 0.002233813827071228
$ws.Range("T4").Value = 0.002233813827071228
$ws.Range("G5").Value = 3.062550666666667
$ws.Range("H5").Value = 9.187652
$ws.Range("I5").Value = 0.06849600470812313
$ws.Range("J5").Value = 0.06849600470812313
$ws.Range("M5").Value = 0.7472223333333332
$ws.Range("N5").Value = 2.241667
$ws.Range("O5").Value = 0.006707722791508481
$ws.Range("P5").Value = 0.006707722791508481
$ws.Range("Q5").Value = 2.288406255098222
$ws.Range("R5").Value = 20.595656295884
$ws.Range("S5").Value = 0.0004594522119079497
$ws.Range("T5").Value = 0.0004594522119079497
$ws.Range("G6").Value = 3.062550666666667
$ws.Range("H6").Value = 9.187652
$ws.Range("I6").Value = 0.06849600470812313
$ws.Range("J6").Value = 0.06849600470812313
$ws.Range("M6").Value = 0.2105813333333333
$ws.Range("N6").Value = 0.631744
$ws.Range("O6").Value = 0.001890362675276361
$ws.Range("P6").Value = 0.001890362675276361
$ws.Range("Q6").Value = 0.6449160027875555
$ws.Range("R6").Value = 5.804244025088
$ws.Range("S6").Value = 0.0001294822907057899
$ws.Range("T6").Value = 0.0001294822907057899
$ws.Range("I7").Value = 0.396815038797359
$ws.Range("J7").Value = 0.396815038797359
$ws.Range("M7").Value = 106.5570066666667
$ws.Range("N7").Value = 319.67102
$ws.Range("O7").Value = 0.956549115742331
$ws.Range("P7").Value = 0.956549115742331
$ws.Range("Q7").Value = 1890.550040640978
$ws.Range("R7").Value = 17014.9503657688
$ws.Range("S7").Value = 0.3795730744748725
$ws.Range("T7").Value = 0.3795730744748725
$ws.Range("I8").Value = 0.396815038797359
$ws.Range("J8").Value = 0.396815038797359
$ws.Range("N8").Value = 0.7487489999999999
$ws.Range("O8").Value = 0.002240475830004717
$ws.Range("P8").Value = 0.002240475830004717
$ws.Range("Q8").Value = 4.428138191506666
$ws.Range("R8").Value = 39.85324372355999
$ws.Range("S8").Value = 0.0008890545034078671
$ws.Range("T8").Value = 0.0008890545034078668
$ws.Range("I9").Value = 0.396815038797359
$ws.Range("J9").Value = 0.396815038797359
$ws.Range("M9").Value = 3.632925333333334
$ws.Range("N9").Value = 10.898776
$ws.Range("O9").Value = 0.03261232296087941
$ws.Range("P9").Value = 0.03261232296087941
$ws.Range("Q9").Value = 64.4558940930489
$ws.Range("R9").Value = 580.10304683744
$ws.Range("S9").Value = 0.01294106020099337
$ws.Range("T9").Value = 0.01294106020099337
$ws.Range("I10").Value = 0.396815038797359
$ws.Range("J10").Value = 0.396815038797359
$ws.Range("M10").Value = 0.7472223333333332
$ws.Range("N10").Value = 2.241667
$ws.Range("O10").Value = 0.006707722791508481
$ws.Range("P10").Value = 0.006707722791508481
$ws.Range("Q10").Value = 13.25732823060889
$ws.Range("R10").Value = 119.31595407548
$ws.Range("S10").Value = 0.002661725279754367
$ws.Range("T10").Value = 0.002661725279754367
$ws.Range("I11").Value = 0.396815038797359
$ws.Range("J11").Value = 0.396815038797359
$ws.Range("M11").Value = 0.2105813333333333
$ws.Range("N11").Value = 0.631744
$ws.Range("O11").Value = 0.001890362675276361
$ws.Range("P11").Value = 0.001890362675276361
$ws.Range("Q11").Value = 3.736164901262222
$ws.Range("R11").Value = 33.62548411136
$ws.Range("S11").Value = 0.0007501243383308687
$ws.Range("T11").Value = 0.0007501243383308687
$ws.Range("G12").Value = 13.27534766666667
$ws.Range("H12").Value = 39.826043
$ws.Range("I12").Value = 0.2969120759943797
$ws.Range("J12").Value = 0.2969120759943796
$ws.Range("M12").Value = 106.5570066666667
$ws.Range("N12").Value = 319.67102
$ws.Range("O12").Value = 0.956549115742331
$ws.Range("P12").Value = 0.956549115742331
$ws.Range("Q12").Value = 1414.581309819318
$ws.Range("R12").Value = 12731.23178837386
$ws.Range("S12").Value = 0.2840109837456437
$ws.Range("T12").Value = 0.2840109837456437
$ws.Range("G13").Value = 13.27534766666667
$ws.Range("H13").Value = 39.826043
$ws.Range("I13").Value = 0.2969120759943797
$ws.Range("J13").Value = 0.2969120759943796
$ws.Range("N13").Value = 0.7487489999999999
$ws.Range("O13").Value = 0.002240475830004717
$ws.Range("P13").Value = 0.002240475830004717
$ws.Range("Q13").Value = 3.313301096689666
$ws.Range("R13").Value = 29.819709870207
$ws.Range("S13").Value = 0.0006652243299019316
$ws.Range("T13").Value = 0.0006652243299019313
$ws.Range("G14").Value = 13.27534766666667
$ws.Range("H14").Value = 39.826043
$ws.Range("I14").Value = 0.2969120759943797
$ws.Range("J14").Value = 0.2969120759943796
$ws.Range("M14").Value = 3.632925333333334
$ws.Range("N14").Value = 10.898776
$ws.Range("O14").Value = 0.03261232296087941
$ws.Range("P14").Value = 0.03261232296087941
$ws.Range("Q14").Value = 48.2283468470409
$ws.Range("R14").Value = 434.0551216233681
$ws.Range("S14").Value = 0.009682992513313881
$ws.Range("T14").Value = 0.009682992513313879
$ws.Range("G15").Value = 13.27534766666667
$ws.Range("H15").Value = 39.826043
$ws.Range("I15").Value = 0.2969120759943797
$ws.Range("J15").Value = 0.2969120759943796
$ws.Range("M15").Value = 0.7472223333333332
$ws.Range("N15").Value = 2.241667
$ws.Range("O15").Value = 0.006707722791508481
$ws.Range("P15").Value = 0.006707722791508481
$ws.Range("Q15").Value = 9.919636259297887
$ws.Range("R15").Value = 89.27672633368098
$ws.Range("S15").Value = 0.001991603899221599
$ws.Range("T15").Value = 0.001991603899221598
$ws.Range("G16").Value = 13.27534766666667
$ws.Range("H16").Value = 39.826043
$ws.Range("I16").Value = 0.2969120759943797
$ws.Range("J16").Value = 0.2969120759943796
$ws.Range("M16").Value = 0.2105813333333333
$ws.Range("N16").Value = 0.631744
$ws.Range("O16").Value = 0.001890362675276361
$ws.Range("P16").Value = 0.001890362675276361
$ws.Range("Q16").Value = 2.795540412110222
$ws.Range("R16").Value = 25.159863708992
$ws.Range("S16").Value = 0.0005612715062985938
$ws.Range("T16").Value = 0.0005612715062985937
$ws.Range("G17").Value = 3.455866
$ws.Range("H17").Value = 10.367598
$ws.Range("I17").Value = 0.07729276657626213
$ws.Range("J17").Value = 0.07729276657626213
$ws.Range("M17").Value = 106.5570066666667
$ws.Range("N17").Value = 319.67102
$ws.Range("O17").Value = 0.956549115742331
$ws.Range("P17").Value = 0.956549115742331
$ws.Range("Q17").Value = 368.2467364011067
$ws.Range("R17").Value = 3314.22062760996
$ws.Range("S17").Value = 0.07393432752180194
$ws.Range("T17").Value = 0.07393432752180194
$ws.Range("G18").Value = 3.455866
$ws.Range("H18").Value = 10.367598
$ws.Range("I18").Value = 0.07729276657626213
$ws.Range("J18").Value = 0.07729276657626213
$ws.Range("N18").Value = 0.7487489999999999
$ws.Range("O18").Value = 0.002240475830004717
$ws.Range("P18").Value = 0.002240475830004717
$ws.Range("Q18").Value = 0.8625254038779999
$ws.Range("R18").Value = 7.762728634901999
$ws.Range("S18").Value = 0.0001731725753483118
$ws.Range("T18").Value = 0.0001731725753483117
$ws.Range("G19").Value = 3.455866
$ws.Range("H19").Value = 10.367598
$ws.Range("I19").Value = 0.07729276657626213
$ws.Range("J19").Value = 0.07729276657626213
$ws.Range("M19").Value = 3.632925333333334
$ws.Range("N19").Value = 10.898776
$ws.Range("O19").Value = 0.03261232296087941
$ws.Range("P19").Value = 0.03261232296087941
$ws.Range("Q19").Value = 12.55490314000533
$ws.Range("R19").Value = 112.994128260048
$ws.Range("S19").Value = 0.002520696666124926
$ws.Range("T19").Value = 0.002520696666124926
$ws.Range("G20").Value = 3.455866
$ws.Range("H20").Value = 10.367598
$ws.Range("I20").Value = 0.07729276657626213
$ws.Range("J20").Value = 0.07729276657626213
$ws.Range("M20").Value = 0.7472223333333332
$ws.Range("N20").Value = 2.241667
$ws.Range("O20").Value = 0.006707722791508481
$ws.Range("P20").Value = 0.006707722791508481
$ws.Range("Q20").Value = 2.582300256207333
$ws.Range("R20").Value = 23.24070230586599
$ws.Range("S20").Value = 0.0005184584519823385
$ws.Range("T20").Value = 0.0005184584519823385
$ws.Range("G21").Value = 3.455866
$ws.Range("H21").Value = 10.367598
$ws.Range("I21").Value = 0.07729276657626213
$ws.Range("J21").Value = 0.07729276657626213
$ws.Range("M21").Value = 0.2105813333333333
$ws.Range("N21").Value = 0.631744
$ws.Range("O21").Value = 0.001890362675276361
$ws.Range("P21").Value = 0.001890362675276361
$ws.Range("Q21").Value = 0.7277408701013333
$ws.Range("R21").Value = 6.549667830911999
$ws.Range("S21").Value = 0.0001461113610046142
$ws.Range("T21").Value = 0.0001461113610046142
$ws.Range("G22").Value = 7.175465666666668
$ws.Range("H22").Value = 21.526397
$ws.Range("I22").Value = 0.1604841139238761
$ws.Range("J22").Value = 0.1604841139238761
$ws.Range("M22").Value = 106.5570066666667
$ws.Range("N22").Value = 319.67102
$ws.Range("O22").Value = 0.956549115742331
$ws.Range("P22").Value = 0.956549115742331
$ws.Range("Q22").Value = 764.5961428794378
$ws.Range("R22").Value = 6881.365285914941
$ws.Range("S22").Value = 0.1535109372645752
$ws.Range("T22").Value = 0.1535109372645752
$ws.Range("G23").Value = 7.175465666666668
$ws.Range("H23").Value = 21.526397
$ws.Range("I23").Value = 0.1604841139238761
$ws.Range("J23").Value = 0.1604841139238761
$ws.Range("N23").Value = 0.7487489999999999
$ws.Range("O23").Value = 0.002240475830004717
$ws.Range("P23").Value = 0.002240475830004717
$ws.Range("Q23").Value = 1.790874247483667
$ws.Range("R23").Value = 16.117868227353
$ws.Range("S23").Value = 0.0003595607783461678
$ws.Range("T23").Value = 0.0003595607783461678
$ws.Range("G24").Value = 7.175465666666668
$ws.Range("H24").Value = 21.526397
$ws.Range("I24").Value = 0.1604841139238761
$ws.Range("J24").Value = 0.1604841139238761
$ws.Range("M24").Value = 3.632925333333334
$ws.Range("N24").Value = 10.898776
$ws.Range("O24").Value = 0.03261232296087941
$ws.Range("P24").Value = 0.03261232296087941
$ws.Range("Q24").Value = 26.0679309988969
$ws.Range("R24").Value = 234.6113789900721
$ws.Range("S24").Value = 0.00523375975337601
$ws.Range("T24").Value = 0.00523375975337601
$ws.Range("G25").Value = 7.175465666666668
$ws.Range("H25").Value = 21.526397
$ws.Range("I25").Value = 0.1604841139238761
$ws.Range("J25").Value = 0.1604841139238761
$ws.Range("M25").Value = 0.7472223333333332
$ws.Range("N25").Value = 2.241667
$ws.Range("O25").Value = 0.006707722791508481
$ws.Range("P25").Value = 0.006707722791508481
$ws.Range("Q25").Value = 5.361668198199888
$ws.Range("R25").Value = 48.255013783799
$ws.Range("S25").Value = 0.001076482948642227
$ws.Range("T25").Value = 0.001076482948642227
$ws.Range("G26").Value = 7.175465666666668
$ws.Range("H26").Value = 21.526397
$ws.Range("I26").Value = 0.1604841139238761
$ws.Range("J26").Value = 0.1604841139238761
$ws.Range("M26").Value = 0.2105813333333333
$ws.Range("N26").Value = 0.631744
$ws.Range("O26").Value = 0.001890362675276361
$ws.Range("P26").Value = 0.001890362675276361
$ws.Range("Q26").Value = 1.511019127374222
$ws.Range("R26").Value = 13.599172146368
$ws.Range("S26").Value = 0.0003033731789364947
$ws.Range("T26").Value = 0.0003033731789364947
